$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 3..11 down to 4..12 (bottom-up so we don't
# overwrite source rows before they are copied), carrying formatting
# (style) first and then the literal values/text.
for ($r = 11; $r -ge 3; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":C" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":C" + $dstRow)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

# Populate the newly-freed row 3 with the inserted test step.
$ws.Range("A3").Value = "verifyElementNotPresent"
$ws.Range("B3").Value = "txt_last_proc_date"
$ws.Range("C3").Value = "getData=WaitForPageLoad"

# Match formatting of the data rows (same style as A4:C4 / old A3:C3)
$ws.Range("A4:C4").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 3 also carries (empty, but styled) D3:E3 cells like the rest of the
# bordered block, plus an F3 cell that only has a border applied.
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("A4:B4").Copy()
$ws.Range("D3:E3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F3").Value = ""
$ws.Range("F3").Borders.LineStyle = 1
